$d = $word.ActiveDocument

# Locate the target paragraph: "(5 points) ... rate a movie on a scale from one to five."
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*rate a movie on a scale from one to five*") {
        # Highlight the whole paragraph, including the paragraph mark, dark green.
        $p.Range.Font.HighlightColorIndex = 11   # wdGreen -> OOXML w:highlight val="darkGreen"
        $found = $true
        break
    }
}

if (-not $found) {
    Write-Host "Target paragraph not found!"
} else {
    Write-Host "Applied dark green highlight to target paragraph."
}
